$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.946.25'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '2.211.08'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '240.81'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").Value = '72.15'
$ws.Range("E7").Value = '  -5.16%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E9").Value = '  -3.87%  '
$ws.Range("D10").Value = '41.76'
$ws.Range("E10").Value = '  -5.43%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.95'
$ws.Range("E12").Value = '  -4.46%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '2.542.39'
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").Value = '14.14'
$ws.Range("E15").Value = '  -3.30%  '
$ws.Range("E16").Value = '  -3.03%  '
$ws.Range("D17").Value = '2.208.93'
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").Value = '41.788.60'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("E19").Value = '  +2.91%  '
$ws.Range("D20").Value = '72.32'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = '6.13'
$ws.Range("E21").Value = '  -1.08%  '
$ws.Range("D22").Value = '10.84'
$ws.Range("E22").Value = '  +18.47%  '
$ws.Range("D23").Value = '228.81'
$ws.Range("E23").Value = '  -1.25%  '
$ws.Range("D24").Value = '2.03'
$ws.Range("E24").Value = '  -8.81%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = '3.63'
$ws.Range("E27").Value = '  +0.46%  '
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("E29").Value = '  -0.88%  '
$ws.Range("D30").Value = '167.43'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").Value = '5.57'
$ws.Range("E32").Value = '  +5.10%  '
$ws.Range("D33").Value = '0.0791'
$ws.Range("E33").Value = '  -4.87%  '
$ws.Range("D34").Value = '29.83'
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("D36").Value = '0.105'
$ws.Range("E36").Value = '  -12.45%  '
$ws.Range("D37").Value = '4.18'
$ws.Range("E37").Value = '  -8.08%  '
$ws.Range("E38").Value = '  -6.50%  '
$ws.Range("D39").Value = '13.76'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("E40").Value = '  -3.85%  '
$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").Value = "'5.60"
$ws.Range("E41").Value = '  -3.60%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = '63.87'
$ws.Range("E42").Value = '  +0.20%  '
$ws.Range("E43").Value = '  -3.94%  '
$ws.Range("D44").Value = '8.63'
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = '102.96'
$ws.Range("E45").Value = '  -4.86%  '
$ws.Range("D46").Value = "'0.1000"
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.33'
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '2.418.09'
$ws.Range("E51").Value = '  -1.65%  '
